$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 241, shifting existing rows 241-354 down to 242-355.
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new record.
$ws.Range("A241").Value = 5
$ws.Range("B241").Value = "Macroferia Regional de Talca"
$ws.Range("C241").Value = "Maule"
$ws.Range("D241").Value = 45029
$ws.Range("E241").Value = 7
$ws.Range("F241").Value = 100112021
$ws.Range("G241").Value = "Ají"
$ws.Range("H241").Value = "Cristal"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 150
$ws.Range("K241").Value = 14000
$ws.Range("L241").Value = 14000
$ws.Range("M241").Value = 14000
$ws.Range("N241").Value = "`$/saco 25 kilos"
$ws.Range("O241").Value = "Región del Maule"
$ws.Range("P241").Value = 560
$ws.Range("Q241").Value = 25
$ws.Range("R241").Value = "Hortaliza"
